# Update the stock-screener table on Sheet1 to the new values (rows expand from 27 to 39)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for A1:F39. $null means the cell is left blank.
$data = @(
    @($null, "Buying Opportunity", "support Zone", "long buildup", "Short buildup", "FII ENTERING"),  # row 1
    @(0, "NSE:AGI", "NSE:5PAISA", "NSE:DRREDDY", "NSE:BHARATFORG", "NSE:M&M"),  # row 2
    @(1, "NSE:BAGFILMS", "NSE:AAVAS", $null, "NSE:HAVELLS", $null),  # row 3
    @(2, "NSE:CHOLAHLDNG", "NSE:ADVANIHOTR", $null, $null, $null),  # row 4
    @(3, "NSE:HCLTECH", "NSE:AVANTIFEED", $null, $null, $null),  # row 5
    @(4, "NSE:HMVL", "NSE:BALAMINES", $null, $null, $null),  # row 6
    @(5, "NSE:INTELLECT", "NSE:BEPL", $null, $null, $null),  # row 7
    @(6, "NSE:M&M", "NSE:BLS", $null, $null, $null),  # row 8
    @(7, "NSE:MAFANG", "NSE:CHEMCON", $null, $null, $null),  # row 9
    @(8, "NSE:MUKTAARTS", "NSE:CHEMFAB", $null, $null, $null),  # row 10
    @(9, "NSE:RAMCOIND", "NSE:DBL", $null, $null, $null),  # row 11
    @(10, "NSE:RHFL", "NSE:DBSTOCKBRO", $null, $null, $null),  # row 12
    @(11, "NSE:SAMBHAAV", "NSE:DECCANCE", $null, $null, $null),  # row 13
    @(12, $null, "NSE:DSSL", $null, $null, $null),  # row 14
    @(13, $null, "NSE:DYCL", $null, $null, $null),  # row 15
    @(14, $null, "NSE:EIHOTEL", $null, $null, $null),  # row 16
    @(15, $null, "NSE:EMSLIMITED", $null, $null, $null),  # row 17
    @(16, $null, "NSE:FINCABLES", $null, $null, $null),  # row 18
    @(17, $null, "NSE:FINEORG", $null, $null, $null),  # row 19
    @(18, $null, "NSE:FUSION", $null, $null, $null),  # row 20
    @(19, $null, "NSE:GODREJAGRO", $null, $null, $null),  # row 21
    @(20, $null, "NSE:GSFC", $null, $null, $null),  # row 22
    @(21, $null, "NSE:HAL", $null, $null, $null),  # row 23
    @(22, $null, "NSE:HERCULES", $null, $null, $null),  # row 24
    @(23, $null, "NSE:INDIAGLYCO", $null, $null, $null),  # row 25
    @(24, $null, "NSE:INOXWIND", $null, $null, $null),  # row 26
    @(25, $null, "NSE:IONEXCHANG", $null, $null, $null),  # row 27
    @(26, $null, "NSE:IRMENERGY", $null, $null, $null),  # row 28
    @(27, $null, "NSE:ITDCEM", $null, $null, $null),  # row 29
    @(28, $null, "NSE:JPOLYINVST", $null, $null, $null),  # row 30
    @(29, $null, "NSE:KIRIINDUS", $null, $null, $null),  # row 31
    @(30, $null, "NSE:MAITHANALL", $null, $null, $null),  # row 32
    @(31, $null, "NSE:MMTC", $null, $null, $null),  # row 33
    @(32, $null, "NSE:OLECTRA", $null, $null, $null),  # row 34
    @(33, $null, "NSE:PATELENG", $null, $null, $null),  # row 35
    @(34, $null, "NSE:PREMEXPLN", $null, $null, $null),  # row 36
    @(35, $null, "NSE:RAMKY", $null, $null, $null),  # row 37
    @(36, $null, "NSE:RTNINDIA", $null, $null, $null),  # row 38
    @(37, $null, "NSE:SALASAR", $null, $null, $null)  # row 39
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 1
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}

# The new rows (28-39) need the same style as the existing numbered rows in column A
# (bold, bordered, centered index cell) - copy that formatting down from row 27.
$ws.Range("A27").Copy()
$ws.Range("A28:A39").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

Write-Host "Updated range:" $ws.Range("A1:F39").Address()